$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update financial data values for rows 2-6 (columns D through AJ)
$ws.Range("D2").Value = 2456
$ws.Range("E2").Value = 219
$ws.Range("F2").Value = 219
$ws.Range("G2").Value = 265
$ws.Range("H2").Value = 213
$ws.Range("I2").Value = 213
$ws.Range("K2").Value = 3304
$ws.Range("L2").Value = 815
$ws.Range("M2").Value = 2489
$ws.Range("N2").Value = 2489
$ws.Range("P2").Value = 110
$ws.Range("Q2").Value = 223
$ws.Range("R2").Value = -146
$ws.Range("S2").Value = -74
$ws.Range("T2").Value = 63
$ws.Range("U2").Value = 160
$ws.Range("V2").Value = 131
$ws.Range("W2").Value = 8.91
$ws.Range("X2").Value = 8.69
$ws.Range("Y2").Value = 8.91
$ws.Range("Z2").Value = 6.6
$ws.Range("AA2").Value = 32.75
$ws.Range("AB2").Value = 1525.96
$ws.Range("AC2").Value = 9700
$ws.Range("AD2").Value = 6.67
$ws.Range("AE2").Value = 113144
$ws.Range("AF2").Value = 0.57
$ws.Range("AG2").Value = 1500
$ws.Range("AH2").Value = 2.32
$ws.Range("AI2").Value = 15.46
$ws.Range("AJ2").Value = 2200000
$ws.Range("D3").Value = 1993
$ws.Range("E3").Value = 111
$ws.Range("F3").Value = 111
$ws.Range("G3").Value = 197
$ws.Range("H3").Value = 160
$ws.Range("I3").Value = 160
$ws.Range("K3").Value = 3407
$ws.Range("L3").Value = 770
$ws.Range("M3").Value = 2637
$ws.Range("N3").Value = 2637
$ws.Range("P3").Value = 110
$ws.Range("Q3").Value = 252
$ws.Range("R3").Value = -217
$ws.Range("S3").Value = 9
$ws.Range("T3").Value = 62
$ws.Range("U3").Value = 190
$ws.Range("V3").Value = 174
$ws.Range("W3").Value = 5.56
$ws.Range("X3").Value = 8.029999999999999
$ws.Range("Y3").Value = 6.25
$ws.Range("Z3").Value = 4.77
$ws.Range("AA3").Value = 29.19
$ws.Range("AB3").Value = 1641.23
$ws.Range("AC3").Value = 7278
$ws.Range("AD3").Value = 7.56
$ws.Range("AE3").Value = 119883
$ws.Range("AF3").Value = 0.46
$ws.Range("AG3").Value = 1500
$ws.Range("AH3").Value = 2.73
$ws.Range("AI3").Value = 20.61
$ws.Range("AJ3").Value = 2200000
$ws.Range("D4").Value = 1748
$ws.Range("E4").Value = 42
$ws.Range("F4").Value = 42
$ws.Range("G4").Value = 83
$ws.Range("H4").Value = 72
$ws.Range("I4").Value = 72
$ws.Range("K4").Value = 3463
$ws.Range("L4").Value = 748
$ws.Range("M4").Value = 2714
$ws.Range("N4").Value = 2714
$ws.Range("P4").Value = 110
$ws.Range("Q4").Value = 184
$ws.Range("R4").Value = -125
$ws.Range("S4").Value = -68
$ws.Range("T4").Value = 38
$ws.Range("U4").Value = 146
$ws.Range("V4").Value = 139
$ws.Range("W4").Value = 2.39
$ws.Range("X4").Value = 4.09
$ws.Range("Y4").Value = 2.67
$ws.Range("Z4").Value = 2.08
$ws.Range("AA4").Value = 27.58
$ws.Range("AB4").Value = 1677.55
$ws.Range("AC4").Value = 3251
$ws.Range("AD4").Value = 15.5
$ws.Range("AE4").Value = 123374
$ws.Range("AF4").Value = 0.41
$ws.Range("AG4").Value = 1250
$ws.Range("AH4").Value = 2.48
$ws.Range("AI4").Value = 38.45
$ws.Range("AJ4").Value = 2200000
$ws.Range("D5").Value = 2084
$ws.Range("E5").Value = 110
$ws.Range("F5").Value = 110
$ws.Range("G5").Value = 113
$ws.Range("H5").Value = 87
$ws.Range("I5").Value = 87
$ws.Range("K5").Value = 3606
$ws.Range("L5").Value = 832
$ws.Range("M5").Value = 2774
$ws.Range("N5").Value = 2774
$ws.Range("P5").Value = 110
$ws.Range("Q5").Value = 61
$ws.Range("R5").Value = -99
$ws.Range("S5").Value = 20
$ws.Range("T5").Value = 27
$ws.Range("U5").Value = 34
$ws.Range("V5").Value = 191
$ws.Range("W5").Value = 5.28
$ws.Range("X5").Value = 4.18
$ws.Range("Y5").Value = 3.18
$ws.Range("Z5").Value = 2.47
$ws.Range("AA5").Value = 30.01
$ws.Range("AB5").Value = 1737.95
$ws.Range("AC5").Value = 3961
$ws.Range("AD5").Value = 14.04
$ws.Range("AE5").Value = 126071
$ws.Range("AF5").Value = 0.44
$ws.Range("AG5").Value = 1250
$ws.Range("AH5").Value = 2.25
$ws.Range("AI5").Value = 31.56
$ws.Range("AJ5").Value = 2200000
$ws.Range("D6").Value = 2156
$ws.Range("E6").Value = 86
$ws.Range("F6").Value = 86
$ws.Range("G6").Value = 154
$ws.Range("H6").Value = 127
$ws.Range("I6").Value = 127
$ws.Range("K6").Value = 3856
$ws.Range("L6").Value = 976
$ws.Range("M6").Value = 2881
$ws.Range("N6").Value = 2881
$ws.Range("P6").Value = 110
$ws.Range("Q6").Value = 90
$ws.Range("R6").Value = -157
$ws.Range("S6").Value = 80
$ws.Range("T6").Value = 105
$ws.Range("U6").Value = -15
$ws.Range("V6").Value = 298
$ws.Range("W6").Value = 3.97
$ws.Range("X6").Value = 5.87
$ws.Range("Y6").Value = 4.48
$ws.Range("Z6").Value = 3.39
$ws.Range("AA6").Value = 33.87
$ws.Range("AB6").Value = 1826.98
$ws.Range("AC6").Value = 5753
$ws.Range("AD6").Value = 7.14
$ws.Range("AE6").Value = 130939
$ws.Range("AF6").Value = 0.31
$ws.Range("AG6").Value = 1250
$ws.Range("AH6").Value = 3.04
$ws.Range("AI6").Value = 21.73
$ws.Range("AJ6").Value = 2200000

# Column J (영업이익(발표기준)) and O (자본총계(지배)) are no longer reported
# for these rows - clear their contents entirely
$ws.Range("J2:J5").ClearContents()
$ws.Range("O2:O5").ClearContents()

# Rows 7, 8 and 9 (2019E/2020E/2021E estimates) no longer carry financial
# data - only the period label in column C remains
$ws.Range("D7:AI7").ClearContents()
$ws.Range("D8:AI8").ClearContents()
$ws.Range("D9:AI9").ClearContents()
